$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# "resultados atualizados" - updated comorbidity summary figures.

# n: 365 -> 364
$t.Cell(2, 2).Range.Text = "364"

# DM (%): 41 (11.2) -> 40 (11.0)
$t.Cell(6, 2).Range.Text = "40 (11.0)"

# HAS (%): 242 (66.3) -> 242 (66.5)
$t.Cell(8, 2).Range.Text = "242 (66.5)"

# Obesidade (%): 115 (34.7) -> 115 (34.8)
$t.Cell(9, 2).Range.Text = "115 (34.8)"
